$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

$ws.Range("E9").Value = "سعید قاسمی"
$ws.Range("E10").Value = "سعید قاسمی"
$ws.Range("E11").Value = "سعید قاسمی"
$ws.Range("E12").Value = "سعید قاسمی"
$ws.Range("E13").Value = "سعید قاسمی"

$ws.Range("D13").Select()
